# "M05 Froze Encoder 1234" - update per-epoch accuracy values in column B
# (column A's epoch index / DisplayOutputs repr cells are untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.84375
$ws.Range("B6").Value = 0.796875
$ws.Range("B7").Value = 0.71875
$ws.Range("B8").Value = 0.734375
$ws.Range("B9").Value = 0.71875
$ws.Range("B10").Value = 0.65625
$ws.Range("B12:B13").Value = 0.65625
$ws.Range("B14").Value = 0.640625
$ws.Range("B15").Value = 0.734375
$ws.Range("B16").Value = 0.625
$ws.Range("B17").Value = 0.65625
$ws.Range("B18").Value = 0.6875
$ws.Range("B20").Value = 0.640625
$ws.Range("B21").Value = 0.625
$ws.Range("B22:B23").Value = 0.65625
$ws.Range("B24").Value = 0.640625
$ws.Range("B25:B36").Value = 0.65625
$ws.Range("B37:B44").Value = 0.671875
$ws.Range("B45:B102").Value = 0.65625
$ws.Range("B105").Value = 0.484375
$ws.Range("B106").Value = 0.46875
$ws.Range("B107").Value = 0.5
$ws.Range("B108").Value = 0.453125
$ws.Range("B109").Value = 0.5625
$ws.Range("B110").Value = 0.546875
$ws.Range("B111").Value = 0.46875
$ws.Range("B112").Value = 0.4375
$ws.Range("B113").Value = 0.578125
$ws.Range("B114").Value = 0.59375
$ws.Range("B115").Value = 0.453125
$ws.Range("B116").Value = 0.546875
$ws.Range("B117").Value = 0.5
$ws.Range("B118").Value = 0.5245901639344263
